# CodeSystem-injury-risk-level-cs.xlsx fix:
#  - Metadata!B7 (the "Experimental" row) gets the literal text value "false"
#    (it was previously blank).
#  - Metadata!B8 (the "Date" row) gets updated to the new publish timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Write "false" as literal text (not boolean TRUE/FALSE) into B7, using a
# leading apostrophe so Excel keeps it as text, then re-apply B7's original
# number format/style (copied from the sibling A7 cell) so no new style slot
# is left attached to the written value.
$ws.Range("B7").Value = "'false"
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the generation timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
